$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "G2"  = 2.7
    "I2"  = 3
    "J2"  = 3.6
    "L2"  = 4
    "M2"  = 1.14
    "N2"  = 5.5
    "Q2"  = 3.1
    "R2"  = 1.36
    "S2"  = 1.73
    "T2"  = 2
    "U2"  = 2.38
    "V2"  = 1.53
    "W2"  = 6
    "X2"  = 11
    "Z2"  = 29
    "AC2" = 5
    "AD2" = 5.5
    "AE2" = 21
    "AI2" = 13
    "AJ2" = 13
    "AK2" = 34
    "AL2" = 34
    "AP2" = 34
    "AQ2" = 67
    "AU2" = 10
    "AW2" = 4.75
    "AZ2" = 67
    "BA2" = 126

    "G3"  = 1.48
    "H3"  = 4.75
    "I3"  = 5.5
    "J3"  = 2
    "L3"  = 5.5
    "O3"  = 1.14
    "P3"  = 5.5
    "Q3"  = 1.5
    "R3"  = 2.5
    "S3"  = 1.25
    "T3"  = 3.75
    "X3"  = 8.5
    "Z3"  = 11
    "AA3" = 11
    "AC3" = 19
    "AD3" = 9
    "AI3" = 34
    "AM3" = 41
    "AO3" = 7
    "AS3" = 81
    "AT3" = 3.75
    "AU3" = 8
    "AW3" = 7.5
    "AX3" = 29
    "BA3" = 101
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
